# Add a new row ("7" | "-" | "0.377") to the end of the Solutions table,
# immediately after the existing row 6 ("probably choose B or C...").

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$newRow = $t.Rows.Add()
$newRow.Cells.Item(1).Range.Text = "7"
$newRow.Cells.Item(2).Range.Text = "-"
$newRow.Cells.Item(3).Range.Text = "0.377"
